# Update NATMI ligand-receptor pair metrics (C1qb-Lrp1) with new TPM-based
# expression values. Every data row (2-21) on the active worksheet gets its
# ligand/receptor/edge expression, specificity, and weight columns (E:T)
# refreshed in place; row/column headers and identifier columns (A:D) are
# untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 1.488897
$ws.Range("H2").Value = 4.466691
$ws.Range("I2").Value = 0.001612421635273429
$ws.Range("J2").Value = 0.001612421635273429
$ws.Range("M2").Value = 6.305846
$ws.Range("N2").Value = 18.917538
$ws.Range("O2").Value = 0.01356150511917599
$ws.Range("P2").Value = 0.01356150511917599
$ws.Range("Q2").Value = 9.388755191862
$ws.Range("R2").Value = 84.498796726758
$ws.Range("S2").Value = 0.00002186686426103072
$ws.Range("T2").Value = 0.00002186686426103072

# Row 3
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 1.488897
$ws.Range("H3").Value = 4.466691
$ws.Range("I3").Value = 0.001612421635273429
$ws.Range("J3").Value = 0.001612421635273429
$ws.Range("O3").Value = 0.392557056479861
$ws.Range("P3").Value = 0.3925570564798609
$ws.Range("Q3").Value = 271.770874231053
$ws.Range("R3").Value = 2445.937868079477
$ws.Range("S3").Value = 0.0006329674909473813
$ws.Range("T3").Value = 0.0006329674909473812

# Row 4
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 1.488897
$ws.Range("H4").Value = 4.466691
$ws.Range("I4").Value = 0.001612421635273429
$ws.Range("J4").Value = 0.001612421635273429
$ws.Range("M4").Value = 127.396393
$ws.Range("N4").Value = 382.189179
$ws.Range("O4").Value = 0.2739817680029065
$ws.Range("P4").Value = 0.2739817680029065
$ws.Range("Q4").Value = 189.680107348521
$ws.Range("R4").Value = 1707.120966136689
$ws.Range("S4").Value = 0.0004417741303983518
$ws.Range("T4").Value = 0.0004417741303983518

# Row 5
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 1.488897
$ws.Range("H5").Value = 4.466691
$ws.Range("I5").Value = 0.001612421635273429
$ws.Range("J5").Value = 0.001612421635273429
$ws.Range("M5").Value = 19.42400133333333
$ws.Range("N5").Value = 58.272004
$ws.Range("O5").Value = 0.04177372766745037
$ws.Range("P5").Value = 0.04177372766745036
$ws.Range("Q5").Value = 28.920337313196
$ws.Range("R5").Value = 260.283035818764
$ws.Range("S5").Value = 0.0000673568622770172
$ws.Range("T5").Value = 0.00006735686227701719

# Row 6
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 1.488897
$ws.Range("H6").Value = 4.466691
$ws.Range("I6").Value = 0.001612421635273429
$ws.Range("J6").Value = 0.001612421635273429
$ws.Range("M6").Value = 129.3233566666667
$ws.Range("N6").Value = 387.97007
$ws.Range("O6").Value = 0.2781259427306063
$ws.Range("P6").Value = 0.2781259427306062
$ws.Range("Q6").Value = 192.54915777093
$ws.Range("R6").Value = 1732.94241993837
$ws.Range("S6").Value = 0.0004484562873896482
$ws.Range("T6").Value = 0.0004484562873896481

# Row 7
$ws.Range("G7").Value = 400.866455
$ws.Range("H7").Value = 1202.599365
$ws.Range("I7").Value = 0.43412388156962
$ws.Range("J7").Value = 0.4341238815696199
$ws.Range("M7").Value = 6.305846
$ws.Range("N7").Value = 18.917538
$ws.Range("O7").Value = 0.01356150511917599
$ws.Range("P7").Value = 0.01356150511917599
$ws.Range("Q7").Value = 2527.80213179593
$ws.Range("R7").Value = 22750.21918616337
$ws.Range("S7").Value = 0.005887373242262951
$ws.Range("T7").Value = 0.00588737324226295

# Row 8
$ws.Range("G8").Value = 400.866455
$ws.Range("H8").Value = 1202.599365
$ws.Range("I8").Value = 0.43412388156962
$ws.Range("J8").Value = 0.4341238815696199
$ws.Range("O8").Value = 0.392557056479861
$ws.Range("P8").Value = 0.3925570564798609
$ws.Range("Q8").Value = 73170.82842214947
$ws.Range("R8").Value = 658537.4557993452
$ws.Range("S8").Value = 0.1704183930965818
$ws.Range("T8").Value = 0.1704183930965817

# Row 9
$ws.Range("G9").Value = 400.866455
$ws.Range("H9").Value = 1202.599365
$ws.Range("I9").Value = 0.43412388156962
$ws.Range("J9").Value = 0.4341238815696199
$ws.Range("M9").Value = 127.396393
$ws.Range("N9").Value = 382.189179
$ws.Range("O9").Value = 0.2739817680029065
$ws.Range("P9").Value = 0.2739817680029065
$ws.Range("Q9").Value = 51068.94044169682
$ws.Range("R9").Value = 459620.4639752714
$ws.Range("S9").Value = 0.1189420286047289
$ws.Range("T9").Value = 0.1189420286047289

# Row 10
$ws.Range("G10").Value = 400.866455
$ws.Range("H10").Value = 1202.599365
$ws.Range("I10").Value = 0.43412388156962
$ws.Range("J10").Value = 0.4341238815696199
$ws.Range("M10").Value = 19.42400133333333
$ws.Range("N10").Value = 58.272004
$ws.Range("O10").Value = 0.04177372766745037
$ws.Range("P10").Value = 0.04177372766745036
$ws.Range("Q10").Value = 7786.430556408607
$ws.Range("R10").Value = 70077.87500767746
$ws.Range("S10").Value = 0.01813497280262578
$ws.Range("T10").Value = 0.01813497280262577

# Row 11
$ws.Range("G11").Value = 400.866455
$ws.Range("H11").Value = 1202.599365
$ws.Range("I11").Value = 0.43412388156962
$ws.Range("J11").Value = 0.4341238815696199
$ws.Range("M11").Value = 129.3233566666667
$ws.Range("N11").Value = 387.97007
$ws.Range("O11").Value = 0.2781259427306063
$ws.Range("P11").Value = 0.2781259427306062
$ws.Range("Q11").Value = 51841.3955356673
$ws.Range("R11").Value = 466572.5598210056
$ws.Range("S11").Value = 0.1207411138234206
$ws.Range("T11").Value = 0.1207411138234206

# Row 12
$ws.Range("I12").Value = 0.0002540413001897126
$ws.Range("J12").Value = 0.0002540413001897126
$ws.Range("M12").Value = 6.305846
$ws.Range("N12").Value = 18.917538
$ws.Range("O12").Value = 0.01356150511917599
$ws.Range("P12").Value = 0.01356150511917599
$ws.Range("Q12").Value = 1.479223252731333
$ws.Range("R12").Value = 13.313009274582
$ws.Range("S12").Value = 0.000003445182393004911
$ws.Range("T12").Value = 0.000003445182393004911

# Row 13
$ws.Range("I13").Value = 0.0002540413001897126
$ws.Range("J13").Value = 0.0002540413001897126
$ws.Range("O13").Value = 0.392557056479861
$ws.Range("P13").Value = 0.3925570564798609
$ws.Range("S13").Value = 0.00009972570502679032
$ws.Range("T13").Value = 0.00009972570502679031

# Row 14
$ws.Range("I14").Value = 0.0002540413001897126
$ws.Range("J14").Value = 0.0002540413001897126
$ws.Range("M14").Value = 127.396393
$ws.Range("N14").Value = 382.189179
$ws.Range("O14").Value = 0.2739817680029065
$ws.Range("P14").Value = 0.2739817680029065
$ws.Range("Q14").Value = 29.88460340447567
$ws.Range("R14").Value = 268.961430640281
$ws.Range("S14").Value = 0.00006960268457173457
$ws.Range("T14").Value = 0.00006960268457173457

# Row 15
$ws.Range("I15").Value = 0.0002540413001897126
$ws.Range("J15").Value = 0.0002540413001897126
$ws.Range("M15").Value = 19.42400133333333
$ws.Range("N15").Value = 58.272004
$ws.Range("O15").Value = 0.04177372766745037
$ws.Range("P15").Value = 0.04177372766745036
$ws.Range("Q15").Value = 4.556475758106222
$ws.Range("R15").Value = 41.008281822956
$ws.Range("S15").Value = 0.00001061225209041006
$ws.Range("T15").Value = 0.00001061225209041006

# Row 16
$ws.Range("I16").Value = 0.0002540413001897126
$ws.Range("J16").Value = 0.0002540413001897126
$ws.Range("M16").Value = 129.3233566666667
$ws.Range("N16").Value = 387.97007
$ws.Range("O16").Value = 0.2781259427306063
$ws.Range("P16").Value = 0.2781259427306062
$ws.Range("Q16").Value = 30.33662989908111
$ws.Range("R16").Value = 273.02966909173
$ws.Range("S16").Value = 0.00007065547610777276
$ws.Range("T16").Value = 0.00007065547610777275

# Row 17
$ws.Range("G17").Value = 520.8019203333333
$ws.Range("H17").Value = 1562.405761
$ws.Range("I17").Value = 0.564009655494917
$ws.Range("J17").Value = 0.5640096554949169
$ws.Range("M17").Value = 6.305846
$ws.Range("N17").Value = 18.917538
$ws.Range("O17").Value = 0.01356150511917599
$ws.Range("P17").Value = 0.01356150511917599
$ws.Range("Q17").Value = 3284.096706126269
$ws.Range("R17").Value = 29556.87035513642
$ws.Range("S17").Value = 0.007648819830259002
$ws.Range("T17").Value = 0.007648819830259

# Row 18
$ws.Range("G18").Value = 520.8019203333333
$ws.Range("H18").Value = 1562.405761
$ws.Range("I18").Value = 0.564009655494917
$ws.Range("J18").Value = 0.5640096554949169
$ws.Range("O18").Value = 0.392557056479861
$ws.Range("P18").Value = 0.3925570564798609
$ws.Range("Q18").Value = 95062.85068087398
$ws.Range("R18").Value = 855565.6561278658
$ws.Range("S18").Value = 0.2214059701873051
$ws.Range("T18").Value = 0.221405970187305

# Row 19
$ws.Range("G19").Value = 520.8019203333333
$ws.Range("H19").Value = 1562.405761
$ws.Range("I19").Value = 0.564009655494917
$ws.Range("J19").Value = 0.5640096554949169
$ws.Range("M19").Value = 127.396393
$ws.Range("N19").Value = 382.189179
$ws.Range("O19").Value = 0.2739817680029065
$ws.Range("P19").Value = 0.2739817680029065
$ws.Range("Q19").Value = 66348.28611794002
$ws.Range("R19").Value = 597134.5750614603
$ws.Range("S19").Value = 0.1545283625832076
$ws.Range("T19").Value = 0.1545283625832075

# Row 20
$ws.Range("G20").Value = 520.8019203333333
$ws.Range("H20").Value = 1562.405761
$ws.Range("I20").Value = 0.564009655494917
$ws.Range("J20").Value = 0.5640096554949169
$ws.Range("M20").Value = 19.42400133333333
$ws.Range("N20").Value = 58.272004
$ws.Range("O20").Value = 0.04177372766745037
$ws.Range("P20").Value = 0.04177372766745036
$ws.Range("Q20").Value = 10116.05719495723
$ws.Range("R20").Value = 91044.51475461504
$ws.Range("S20").Value = 0.02356078575045717
$ws.Range("T20").Value = 0.02356078575045716

# Row 21
$ws.Range("G21").Value = 520.8019203333333
$ws.Range("H21").Value = 1562.405761
$ws.Range("I21").Value = 0.564009655494917
$ws.Range("J21").Value = 0.5640096554949169
$ws.Range("M21").Value = 129.3233566666667
$ws.Range("N21").Value = 387.97007
$ws.Range("O21").Value = 0.2781259427306063
$ws.Range("P21").Value = 0.2781259427306062
$ws.Range("Q21").Value = 67351.8524959526
$ws.Range("R21").Value = 606166.6724635733
$ws.Range("S21").Value = 0.1568657171436882
$ws.Range("T21").Value = 0.1568657171436882
